# Auto-generated Excel COM-interop script
# Refreshes market-price derived figures on the Cerberus_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match latest scheduled-runner pull,
# per commit 'chore: update Sheets via scheduled runner'.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H15").Value = 1414.3881
$ws.Range("I15").Value = 1414.3881
$ws.Range("K15").Value = 4243.164299999999
$ws.Range("M15").Value = -4074.164299999999
$ws.Range("H16").Value = 16499.5
$ws.Range("J16").Value = 16499.5
$ws.Range("L16").Value = 16499.5
$ws.Range("N16").Value = -16959.5
$ws.Range("H17").Value = 2843.8572
$ws.Range("J17").Value = 2843.8572
$ws.Range("L17").Value = 8531.571599999999
$ws.Range("N17").Value = -8867.571599999999
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 406.8
$ws.Range("J33").Value = 888
$ws.Range("L33").Value = 888
$ws.Range("N33").Value = -1346
$ws.Range("H51").Value = 7076.467
$ws.Range("I51").Value = 8600
$ws.Range("J51").Value = 6314.7
$ws.Range("K51").Value = 8600
$ws.Range("L51").Value = 6314.7
$ws.Range("M51").Value = -8116
$ws.Range("N51").Value = -7282.7
$ws.Range("H69").Value = 14274.857
$ws.Range("I69").Value = 14025.667
$ws.Range("J69").Value = 14461.75
$ws.Range("K69").Value = 42077.001
$ws.Range("L69").Value = 43385.25
$ws.Range("M69").Value = -41203.001
$ws.Range("N69").Value = -45133.25
$ws.Range("H72").Value = 14274.857
$ws.Range("I72").Value = 14025.667
$ws.Range("J72").Value = 14461.75
$ws.Range("K72").Value = 126231.003
$ws.Range("L72").Value = 130155.75
$ws.Range("M72").Value = -121863.003
$ws.Range("N72").Value = -138891.75
$ws.Range("H80").Value = 593.2941
$ws.Range("I80").Value = 416.66666
$ws.Range("J80").Value = 631.1429000000001
$ws.Range("K80").Value = 1249.99998
$ws.Range("L80").Value = 1893.4287
$ws.Range("M80").Value = -251.9999800000001
$ws.Range("N80").Value = -3889.4287
$ws.Range("H83").Value = 593.2941
$ws.Range("I83").Value = 416.66666
$ws.Range("J83").Value = 631.1429000000001
$ws.Range("K83").Value = 3749.99994
$ws.Range("L83").Value = 5680.2861
$ws.Range("M83").Value = 1242.00006
$ws.Range("N83").Value = -15664.2861
$ws.Range("H99").Value = 1734.7
$ws.Range("J99").Value = 2080.2856
$ws.Range("L99").Value = 6240.8568
$ws.Range("N99").Value = -9236.856800000001
$ws.Range("H127").Value = 1534
$ws.Range("I127").Value = 1540.8
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 4622.4
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = 337.6000000000004
$ws.Range("N127").Value = -14420
$ws.Range("H132").Value = 3595.74
$ws.Range("I132").Value = 3546.6938
$ws.Range("K132").Value = 10640.0814
$ws.Range("M132").Value = -8110.081399999999
$ws.Range("H134").Value = 27160.412
$ws.Range("J134").Value = 27160.412
$ws.Range("L134").Value = 27160.412
$ws.Range("N134").Value = -37300.412
$ws.Range("H138").Value = 4296.7856
$ws.Range("J138").Value = 4236.609
$ws.Range("L138").Value = 12709.827
$ws.Range("N138").Value = -22989.827

# ---- Sheet: ARM ----
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 3693.9768
$ws.Range("I32").Value = 2316.342
$ws.Range("K32").Value = 2316.342
$ws.Range("M32").Value = -2029.342
$ws.Range("H61").Value = 3333.7585
$ws.Range("I61").Value = 2967.24
$ws.Range("K61").Value = 2967.24
$ws.Range("M61").Value = -2755.24
$ws.Range("H102").Value = 7982.6665
$ws.Range("I102").Value = 7982.6665
$ws.Range("K102").Value = 7982.6665
$ws.Range("M102").Value = -6360.6665
$ws.Range("H136").Value = 3333.7585
$ws.Range("I136").Value = 2967.24
$ws.Range("K136").Value = 8901.719999999999
$ws.Range("M136").Value = -6351.719999999999

# ---- Sheet: BSM ----
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H5").Value = 343.33334
$ws.Range("I5").Value = 343.33334
$ws.Range("K5").Value = 343.33334
$ws.Range("M5").Value = -230.33334
$ws.Range("H22").Value = 616.8333
$ws.Range("J22").Value = 150
$ws.Range("L22").Value = 150
$ws.Range("N22").Value = -496
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H99").Value = 1362.2693
$ws.Range("J99").Value = 1450
$ws.Range("L99").Value = 1450
$ws.Range("N99").Value = -4446
$ws.Range("H105").Value = 3074.5
$ws.Range("I105").Value = 2641.25
$ws.Range("K105").Value = 2641.25
$ws.Range("M105").Value = -894.25
$ws.Range("H109").Value = 45053.6
$ws.Range("J109").Value = 46067
$ws.Range("L109").Value = 46067
$ws.Range("N109").Value = -48147
$ws.Range("H126").Value = 1362.2693
$ws.Range("J126").Value = 1450
$ws.Range("L126").Value = 4350
$ws.Range("N126").Value = -9290

# ---- Sheet: CUL ----
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H80").Value = 3874.5
$ws.Range("J80").Value = 4499.3335
$ws.Range("L80").Value = 13498.0005
$ws.Range("N80").Value = -15370.0005
$ws.Range("H83").Value = 3874.5
$ws.Range("J83").Value = 4499.3335
$ws.Range("L83").Value = 40494.0015
$ws.Range("N83").Value = -49854.0015
$ws.Range("H122").Value = 1652.875
$ws.Range("J122").Value = 2660.5557
$ws.Range("L122").Value = 23945.0013
$ws.Range("N122").Value = -28845.0013

# ---- Sheet: GSM ----
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H102").Value = 10955.258
$ws.Range("I102").Value = 13325.728
$ws.Range("K102").Value = 13325.728
$ws.Range("M102").Value = -11703.728

# ---- Sheet: LTW ----
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H102").Value = 69706.664
$ws.Range("J102").Value = 69706.664
$ws.Range("L102").Value = 69706.664
$ws.Range("N102").Value = -76196.664
$ws.Range("H109").Value = 64989.5
$ws.Range("J109").Value = 64989.5
$ws.Range("L109").Value = 64989.5
$ws.Range("N109").Value = -67763.5

# ---- Sheet: WVR ----
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H8").Value = 12499.75
$ws.Range("I8").Value = 5000
$ws.Range("K8").Value = 5000
$ws.Range("M8").Value = -4860
$ws.Range("H70").Value = 44525.5
$ws.Range("I70").Value = 45998
$ws.Range("J70").Value = 44034.668
$ws.Range("K70").Value = 45998
$ws.Range("L70").Value = 44034.668
$ws.Range("M70").Value = -45683
$ws.Range("N70").Value = -44664.668
$ws.Range("H73").Value = 44525.5
$ws.Range("I73").Value = 45998
$ws.Range("J73").Value = 44034.668
$ws.Range("K73").Value = 45998
$ws.Range("L73").Value = 44034.668
$ws.Range("M73").Value = -44906
$ws.Range("N73").Value = -46218.668
$ws.Range("H109").Value = 79997.5
$ws.Range("J109").Value = 79997.5
$ws.Range("L109").Value = 79997.5
$ws.Range("N109").Value = -82771.5
$ws.Range("H132").Value = 2076.7058
$ws.Range("I132").Value = 1807.75
$ws.Range("K132").Value = 5423.25
$ws.Range("M132").Value = -2893.25
